$wb = $excel.ActiveWorkbook

# "SoCDTtiNTY-psgr" sheet (2nd sheet) - revert data values and formula
$ws = $wb.Worksheets.Item(2)

# B2: restore formula that computes 0.076 + (0.076 - 0.0725) = 0.0795
$ws.Range("B2").Formula = "=0.076+(0.076-0.0725)"

# D2: restore value to 0.0735
$ws.Range("D2").Value = 0.0735

# B5 and E5: restore value to 0.01
$ws.Range("B5").Value = 0.01
$ws.Range("E5").Value = 0.01

# Make the psgr sheet the active/selected tab, with E6 as the selected cell
$ws.Activate() | Out-Null
$ws.Range("E6").Select() | Out-Null
